# localSearch workbook edit:
#  - Column F ("BestKnownValue") was stored as text (shared strings) and is
#    converted to real numbers.
#  - A new column H ("Local Search besser als bester bekannter Wert?") is
#    added with a Ja/Nein formula comparing LocalSearchValues/VehicleNumber
#    against the (now numeric) BestKnownValue/VehicleNumber.
#  - Selection moves to H8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F: replace the text "numbers" with real numeric values -------
$bestKnownValues = @{
    2  = 828.94;    3  = 828.94;    4  = 828.94;    5  = 824.78;   6  = 828.94
    7  = 828.94;    8  = 828.94;    9  = 828.94;    10 = 828.94;   11 = 828.94
    12 = 591.56;    13 = 591.56;    14 = 591.17;    15 = 590.6;    16 = 588.29
    17 = 588.49;    18 = 588.29;    19 = 588.32;    20 = 1650.8;   21 = 1486.12
    22 = 1292.68;   23 = 1007.31;   24 = 1377.11;   25 = 1252.03;  26 = 1104.66
    27 = 960.88;    28 = 1194.73;   29 = 1118.84;   30 = 1096.72;  31 = 982.14
    32 = 1252.37;   33 = 1191.7;    34 = 939.5;     35 = 825.52;   36 = 994.42
    37 = 906.14;    38 = 890.61;    39 = 726.82;    40 = 909.16;   41 = 939.37
    42 = 885.71;    43 = 1696.94;   44 = 1554.75;   45 = 1261.67;  46 = 1135.48
    47 = 1629.44;   48 = 1424.73;   49 = 1230.48;   50 = 1139.82;  51 = 1406.94
    52 = 1365.65;   53 = 1049.62;   54 = 798.46;    55 = 1297.65;  56 = 1146.32
    57 = 1061.14
}

for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 6).Value = $bestKnownValues[$row]
}

# --- Column H: header + Ja/Nein formula -----------------------------------
$ws.Range("H1").Value = "Local Search besser als bester bekannter Wert?"

# H2 is entered as a standalone formula ...
$ws.Range("H2").Formula = '=IF(D2<F2,IF(E2<=G2,"Ja","Nein"),"Nein")'
# ... then filled down through H57, which Excel stores as a shared formula.
$ws.Range("H3:H57").Formula = '=IF(D3<F3,IF(E3<=G3,"Ja","Nein"),"Nein")'

# --- Selection -------------------------------------------------------------
[void]$ws.Range("H8").Select()
